# The commit inserts a new weekly price record for "Cilantro" (Femacal de
# La Calera) as row 158, shifting all the following records down by one
# row (old row 158 becomes 159, ..., old row 225 becomes 226).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 158; Excel shifts rows 158:225 down to 159:226
# and copies formatting (e.g. the date style on column D) from the row above.
$ws.Rows.Item(158).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(158, 1).Value  = 3
$ws.Cells.Item(158, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(158, 3).Value  = "Coquimbo"
$ws.Cells.Item(158, 4).Value  = 44489
$ws.Cells.Item(158, 5).Value  = 5
$ws.Cells.Item(158, 6).Value  = 100112040
$ws.Cells.Item(158, 7).Value  = "Cilantro"
$ws.Cells.Item(158, 8).Value  = "Sin especificar"
$ws.Cells.Item(158, 9).Value  = "Primera"
$ws.Cells.Item(158, 10).Value = 230
$ws.Cells.Item(158, 11).Value = 2500
$ws.Cells.Item(158, 12).Value = 2800
$ws.Cells.Item(158, 13).Value = 2657
$ws.Cells.Item(158, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(158, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(158, 16).Value = 886
$ws.Cells.Item(158, 17).Value = 3
$ws.Cells.Item(158, 18).Value = "Hortaliza"
